$d = $word.ActiveDocument

# Locate the run that reads "učitelja TZK-e, " (style s2) which is being
# replaced by the template placeholder "{{ r_mj_zamijenj_G }}, ".
$find = $d.Content
$found = $find.Find.Execute("učitelja TZK-e, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $find.Start

    # New combined text that will occupy the same slot.
    $newText = "{{ r_mj_zamijenj_G }}, "

    # Replace the whole old run's text in one shot (keeps it anchored at $start).
    $find.Text = $newText

    # Piece boundaries (character offsets relative to document start):
    #   A = "{{ "                 -> style apple-converted-space
    #   B = "r_mj_zamijenj_G }}"  -> style s2
    #   C = ", "                  -> style s2
    $aText = "{{ "
    $bText = "r_mj_zamijenj_G }}"
    $cText = ", "

    $aStart = $start
    $aEnd   = $aStart + $aText.Length

    $bStart = $aEnd
    $bEnd   = $bStart + $bText.Length

    $cStart = $bEnd
    $cEnd   = $cStart + $cText.Length

    $rA = $d.Range($aStart, $aEnd)
    $rA.Style = "apple-converted-space"

    $rB = $d.Range($bStart, $bEnd)
    $rB.Style = "s2"

    $rC = $d.Range($cStart, $cEnd)
    $rC.Style = "s2"
}
